$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 2484.5
$ws.Range("I49").Value = 1969
$ws.Range("J49").Value = 3000
$ws.Range("K49").Value = 5907
$ws.Range("L49").Value = 9000
$ws.Range("M49").Value = -5771
$ws.Range("N49").Value = -9272

$ws.Range("H59").Value = 1000
$ws.Range("J59").Value = 1000
$ws.Range("L59").Value = 3000
$ws.Range("N59").Value = -4114

$ws.Range("H116").Value = 2514.9443
$ws.Range("I116").Value = 2740.4443
$ws.Range("J116").Value = 2289.4443
$ws.Range("K116").Value = 2740.4443
$ws.Range("L116").Value = 2289.4443
$ws.Range("M116").Value = 701.5556999999999
$ws.Range("N116").Value = -9173.444299999999

$ws.Range("H138").Value = 13891689
$ws.Range("J138").Value = 3111.2808
$ws.Range("L138").Value = 9333.8424
$ws.Range("N138").Value = -19613.8424

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 52719.668
$ws.Range("J24").Value = 52719.668
$ws.Range("L24").Value = 52719.668
$ws.Range("N24").Value = -53467.668

$ws.Range("H74").Value = 3154.6614
$ws.Range("I74").Value = 2870.7576
$ws.Range("J74").Value = 3477.724
$ws.Range("K74").Value = 2870.7576
$ws.Range("L74").Value = 3477.724
$ws.Range("M74").Value = -1996.7576
$ws.Range("N74").Value = -5225.724

$ws.Range("H77").Value = 3154.6614
$ws.Range("I77").Value = 2870.7576
$ws.Range("J77").Value = 3477.724
$ws.Range("K77").Value = 14353.788
$ws.Range("L77").Value = 17388.62
$ws.Range("M77").Value = -9985.788
$ws.Range("N77").Value = -26124.62

$ws.Range("H100").Value = 52719.668
$ws.Range("J100").Value = 52719.668
$ws.Range("L100").Value = 52719.668
$ws.Range("N100").Value = -54883.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1127.7556
$ws.Range("I20").Value = 1231.875
$ws.Range("J20").Value = 1008.7619
$ws.Range("K20").Value = 1231.875
$ws.Range("L20").Value = 1008.7619
$ws.Range("M20").Value = -984.875
$ws.Range("N20").Value = -1502.7619

$ws.Range("H86").Value = 3703.1
$ws.Range("I86").Value = 3047.6
$ws.Range("K86").Value = 3047.6
$ws.Range("M86").Value = -1924.6

$ws.Range("H89").Value = 3703.1
$ws.Range("I89").Value = 3047.6
$ws.Range("K89").Value = 15238
$ws.Range("M89").Value = -9622

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5698.533
$ws.Range("I16").Value = 5698.533
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 5698.533
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -5411.533
$ws.Range("N16").ClearContents()

$ws.Range("H43").Value = 110000
$ws.Range("J43").Value = 110000
$ws.Range("L43").Value = 110000
$ws.Range("N43").Value = -110368

$ws.Range("H86").Value = 66671932
$ws.Range("I86").Value = 125003800
$ws.Range("K86").Value = 125003800
$ws.Range("M86").Value = -125002677

$ws.Range("H89").Value = 66671932
$ws.Range("I89").Value = 125003800
$ws.Range("K89").Value = 625019000
$ws.Range("M89").Value = -625013384

$ws.Range("H101").Value = 110000
$ws.Range("J101").Value = 110000
$ws.Range("L101").Value = 110000
$ws.Range("N101").Value = -116490

$ws.Range("H107").Value = 880.2632
$ws.Range("J107").Value = 1062.5555
$ws.Range("L107").Value = 1062.5555
$ws.Range("N107").Value = -4902.5555

$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("M108").ClearContents()
$ws.Range("N108").ClearContents()

$ws.Range("H109").Value = 95000
$ws.Range("J109").Value = 95000
$ws.Range("L109").Value = 95000
$ws.Range("N109").Value = -97080

$ws.Range("H112").Value = 92990
$ws.Range("J112").Value = 92990
$ws.Range("L112").Value = 92990
$ws.Range("N112").Value = -95944

$ws.Range("H113").Value = 5698.533
$ws.Range("I113").Value = 5698.533
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 5698.533
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -3528.533
$ws.Range("N113").ClearContents()

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H132").Value = 1356.7916
$ws.Range("I132").Value = 1235.8572
$ws.Range("K132").Value = 3707.5716
$ws.Range("M132").Value = -1177.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1118.6
$ws.Range("I5").Value = 968.38464
$ws.Range("J5").Value = 2095
$ws.Range("K5").Value = 2905.15392
$ws.Range("L5").Value = 6285
$ws.Range("M5").Value = -2793.15392
$ws.Range("N5").Value = -6509

$ws.Range("H135").Value = 1118.6
$ws.Range("I135").Value = 968.38464
$ws.Range("J135").Value = 2095
$ws.Range("K135").Value = 8715.46176
$ws.Range("L135").Value = 18855
$ws.Range("M135").Value = -6180.46176
$ws.Range("N135").Value = -23925

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 55643.617
$ws.Range("I70").Value = 143342.75
$ws.Range("J70").Value = 16666.223
$ws.Range("K70").Value = 143342.75
$ws.Range("L70").Value = 16666.223
$ws.Range("M70").Value = -143072.75
$ws.Range("N70").Value = -17206.223

$ws.Range("H73").Value = 55643.617
$ws.Range("I73").Value = 143342.75
$ws.Range("J73").Value = 16666.223
$ws.Range("K73").Value = 143342.75
$ws.Range("L73").Value = 16666.223
$ws.Range("M73").Value = -142406.75
$ws.Range("N73").Value = -18538.223

$ws.Range("H132").Value = 2584.8293
$ws.Range("I132").Value = 2582.05
$ws.Range("K132").Value = 7746.150000000001
$ws.Range("M132").Value = -5216.150000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2837.26
$ws.Range("I132").Value = 2915.284
$ws.Range("K132").Value = 8745.852000000001
$ws.Range("M132").Value = -6215.852000000001

$ws.Range("H136").Value = 3329.6775
$ws.Range("I136").Value = 2993.7932
$ws.Range("J136").Value = 8200
$ws.Range("K136").Value = 8981.3796
$ws.Range("L136").Value = 24600
$ws.Range("M136").Value = -6431.3796
$ws.Range("N136").Value = -29700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3591.4565
$ws.Range("J136").Value = 2957.2
$ws.Range("L136").Value = 8871.599999999999
$ws.Range("N136").Value = -13971.6
